# The underlying dataset rows (2..55) were shuffled/re-paired: every data
# row's content moved to a different row position. Only a subset of
# columns actually differ from row to row (the rest are constant across
# the whole sheet), so we snapshot those columns, remap them according to
# the permutation below, and write the values back.
#
# $perm[i] is the ORIGINAL row number whose data now belongs at sheet row
# (i + 2)  -- i.e. $perm[0] holds the source row for row 2, $perm[1] the
# source row for row 3, and so on through row 55.
$perm = @(
    45,10,29, 7,37,13,17,22,33,41,
    36,20,18,11,51,26, 4,31,19,16,
    24,32,12,40, 9,42,35,47,27,48,
     6,52,44,43,14, 5,30,49,28,55,
     3,54,21,34, 2,46,39, 8,50,25,
    53,15,23,38
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually vary row-to-row for this dataset. Every
# other column in A1:AY55 is identical on every data row, so it is left
# untouched (this also sidesteps COM's auto date/time coercion on the
# Y/AA "Startdatum"/"Slutdatum" text columns, which are constant anyway).
$cols = @("A","B","D","E","F","G","H","I","M","Q","R","Z","AB","AC","AQ","AR")

$firstRow = 2
$lastRow = 55

# Snapshot current values before any writes (rows must not be overwritten
# before they are read, since we are permuting in place).
$snapshot = @{}
foreach ($col in $cols) {
    $colVals = @{}
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $colVals[$r] = $ws.Range($col + $r).Value2
    }
    $snapshot[$col] = $colVals
}

# Write back in permuted order.
for ($i = 0; $i -lt $perm.Length; $i++) {
    $destRow = $firstRow + $i
    $srcRow = $perm[$i]
    foreach ($col in $cols) {
        $ws.Range($col + $destRow).Value2 = $snapshot[$col][$srcRow]
    }
}
